$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '69.420.95'
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '3.884.15'
$ws.Range('E3').NumberFormat = "@"
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').NumberFormat = "@"
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '604.14'
$ws.Range('E5').NumberFormat = "@"
$ws.Range('E5').Value = '  +0.65%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '170.35'
$ws.Range('E6').NumberFormat = "@"
$ws.Range('E6').Value = '  +4.29%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '3.883.66'
$ws.Range('E7').NumberFormat = "@"
$ws.Range('E7').Value = '  +1.47%  '
$ws.Range('E9').NumberFormat = "@"
$ws.Range('E9').Value = '  +1.08%  '
$ws.Range('E10').NumberFormat = "@"
$ws.Range('E10').Value = '  +1.61%  '
$ws.Range('E11').NumberFormat = "@"
$ws.Range('E11').Value = '  +1.14%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.468'
$ws.Range('E12').NumberFormat = "@"
$ws.Range('E12').Value = '  +1.99%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.0000255'
$ws.Range('E13').NumberFormat = "@"
$ws.Range('E13').Value = '  +4.72%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '38.26'
$ws.Range('E14').NumberFormat = "@"
$ws.Range('E14').Value = '  +4.02%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '4.540.83'
$ws.Range('E15').NumberFormat = "@"
$ws.Range('E15').Value = '  +1.51%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.892.62'
$ws.Range('E16').NumberFormat = "@"
$ws.Range('E16').Value = '  +1.58%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '69.493.06'
$ws.Range('E17').NumberFormat = "@"
$ws.Range('E17').Value = '  +1.26%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '18.77'
$ws.Range('E18').NumberFormat = "@"
$ws.Range('E18').Value = '  +9.74%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '7.61'
$ws.Range('E19').NumberFormat = "@"
$ws.Range('E19').Value = '  +0.55%  '
$ws.Range('E20').NumberFormat = "@"
$ws.Range('E20').Value = '  -0.67%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.05'
$ws.Range('E21').NumberFormat = "@"
$ws.Range('E21').Value = '  -1.13%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '489.19'
$ws.Range('E22').NumberFormat = "@"
$ws.Range('E22').Value = '  +0.82%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '0.742'
$ws.Range('E23').NumberFormat = "@"
$ws.Range('E23').Value = '  +3.58%  '
$ws.Range('E24').NumberFormat = "@"
$ws.Range('E24').Value = '  +4.09%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '85.31'
$ws.Range('E25').NumberFormat = "@"
$ws.Range('E25').Value = '  +1.50%  '
$ws.Range('E26').NumberFormat = "@"
$ws.Range('E26').Value = '  +2.62%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.38'
$ws.Range('E27').NumberFormat = "@"
$ws.Range('E27').Value = '  +2.36%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '10.12'
$ws.Range('E28').NumberFormat = "@"
$ws.Range('E28').Value = '  +1.27%  '
$ws.Range('E29').NumberFormat = "@"
$ws.Range('E29').Value = '  +0.25%  '
$ws.Range('E30').NumberFormat = "@"
$ws.Range('E30').Value = '  +1.22%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '4.037.41'
$ws.Range('E31').NumberFormat = "@"
$ws.Range('E31').Value = '  +1.39%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '2.39'
$ws.Range('E32').NumberFormat = "@"
$ws.Range('E32').Value = '  +1.48%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '7.82'
$ws.Range('E33').NumberFormat = "@"
$ws.Range('E33').Value = '  -0.16%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '31.86'
$ws.Range('E34').NumberFormat = "@"
$ws.Range('E34').Value = '  +0.26%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '3.855.40'
$ws.Range('E35').NumberFormat = "@"
$ws.Range('E35').Value = '  +2.11%  '
$ws.Range('E36').NumberFormat = "@"
$ws.Range('E36').Value = '  -0.13%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '3.41'
$ws.Range('E37').NumberFormat = "@"
$ws.Range('E37').Value = '  +15.33%  '
$ws.Range('B38').NumberFormat = "@"
$ws.Range('B38').Value = 'Kaspa'
$ws.Range('C38').NumberFormat = "@"
$ws.Range('C38').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.143'
$ws.Range('E38').NumberFormat = "@"
$ws.Range('E38').Value = '  +2.66%  '
$ws.Range('B39').NumberFormat = "@"
$ws.Range('B39').Value = 'Filecoin'
$ws.Range('C39').NumberFormat = "@"
$ws.Range('C39').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '6.11'
$ws.Range('E39').NumberFormat = "@"
$ws.Range('E39').Value = '  +4.07%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.03'
$ws.Range('E40').NumberFormat = "@"
$ws.Range('E40').Value = '  +0.58%  '
$ws.Range('E41').NumberFormat = "@"
$ws.Range('E41').Value = '  +0.00%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.326'
$ws.Range('E42').NumberFormat = "@"
$ws.Range('E42').Value = '  +2.53%  '
$ws.Range('E43').NumberFormat = "@"
$ws.Range('E43').Value = '  +4.69%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '435.47'
$ws.Range('E44').NumberFormat = "@"
$ws.Range('E44').Value = '  +1.65%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '47.97'
$ws.Range('E45').NumberFormat = "@"
$ws.Range('E45').Value = '  -0.98%  '
$ws.Range('E46').NumberFormat = "@"
$ws.Range('E46').Value = '  +3.42%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.000276'
$ws.Range('E48').NumberFormat = "@"
$ws.Range('E48').Value = '  +22.06%  '
$ws.Range('E49').NumberFormat = "@"
$ws.Range('E49').Value = '  +2.43%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '40.25'
$ws.Range('E50').NumberFormat = "@"
$ws.Range('E50').Value = '  +4.12%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '141.13'
$ws.Range('E51').NumberFormat = "@"
$ws.Range('E51').Value = '  -1.07%  '
